$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 5 (pushes the NAME/CONDITION/ACTION table down)
$ws.Rows("5:5").Insert()

# 2. New title cell A5:C5 = "RuleTable 黑名單"
$ws.Range("A5:C5").Merge()
$ws.Range("A5").Value = "RuleTable 黑名單"

# 3. Column widths for A and D (new)
$ws.Columns("A:A").ColumnWidth = 14.571428571428573
$ws.Columns("D:D").ColumnWidth = 19.428571428571427

# 4. Select A9:C9 with C9 active (post-insert the header row "Rule Name/ID/Status" is now row 9)
$ws.Range("A9:C9").Select()
